# "Generate Report for Archive"
#
# The localization-status report moved from "Ready for handoff" to
# "In Translation" for the 35fcf230...md entry. That status string is
# shared across the Overview sheet (zh-cn / de-de status columns) and the
# per-language detail sheets (zh-cn, de-de), so update every cell that
# shows it. Excel then re-ran its column "best fit" for the now-shorter
# text, narrowing the affected Status columns.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview: zh-cn / de-de status columns (E, F) for the single data row.
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# Per-language detail sheets: Status column (C) for the single data row.
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# Re-fit the narrowed Status columns to the new text.
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
